$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.138.90'
$ws.Range('E2').Value = '  +4.39%  '
$ws.Range('D3').Value = '2.502.35'
$ws.Range('E3').Value = '  +2.40%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '494.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.84%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.516'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.22%  '
$ws.Range('D9').Value = '2.521.48'
$ws.Range('E9').Value = '  +2.72%  '
$ws.Range('E10').Value = '  +5.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.76'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.338'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.48%  '
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').Value = '2.939.00'
$ws.Range('E14').Value = '  +2.03%  '
$ws.Range('D15').Value = '57.291.79'
$ws.Range('E15').Value = '  +3.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.98%  '
$ws.Range('E17').Value = '  +2.99%  '
$ws.Range('D18').Value = '2.509.94'
$ws.Range('E18').Value = '  +2.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.58'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.38%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.95'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.52'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.410'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.76%  '
$ws.Range('E26').Value = '  +3.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').Value = '2.624.15'
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.17%  '
$ws.Range('D30').Value = '0.0₃0825'
$ws.Range('E30').Value = '  +7.39%  '
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.42'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.37'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.53'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.28'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.908'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.97%  '
$ws.Range('E38').Value = '  +4.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '34.35'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.53'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.05%  '
$ws.Range('E42').Value = '  +2.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0561'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.995'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '266.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0948'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0230'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.22'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.02'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.53%  '
$ws.Range('D51').Value = '1.895.47'
$ws.Range('E51').Value = '  -1.62%  '
